$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.120882630348206
$ws.Range("B1").Value = 2.341852188110352
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.824825048446655
$ws.Range("E1").Value = 1.224672079086304
